$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$refStyle = $ws.Range("D7").Style

$ws.Range("D2").Value = "68.971.14"
$ws.Range("E2").Value = "  +1.90%  "

$ws.Range("D3").Value = "2.522.06"
$ws.Range("E3").Value = "  +1.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $refStyle
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.97"
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = "  +1.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.91"
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("D8").Style = $refStyle
$ws.Range("E8").Value = "  +0.58%  "

$ws.Range("D9").Value = "2.521.32"
$ws.Range("E9").Value = "  +1.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("D10").Style = $refStyle
$ws.Range("E10").Value = "  +5.84%  "

$ws.Range("E11").Value = "  -0.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.98"
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.980.48"
$ws.Range("E14").Value = "  +1.04%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.93"
$ws.Range("D15").Style = $refStyle
$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("D16").Value = "68.796.83"
$ws.Range("E16").Value = "  +1.80%  "

$ws.Range("E17").Value = "  +1.03%  "

$ws.Range("D18").Value = "2.475.82"
$ws.Range("E18").Value = "  +0.14%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "362.95"
$ws.Range("D19").Style = $refStyle
$ws.Range("E19").Value = "  +3.18%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.05"
$ws.Range("D20").Style = $refStyle
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("E21").Value = "  +1.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.08"
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.64"
$ws.Range("D24").Style = $refStyle
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.20"
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = "  -0.68%  "

$ws.Range("E26").Value = "  -5.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("D27").Style = $refStyle
$ws.Range("E27").Value = "  -1.83%  "

$ws.Range("D28").Value = "2.646.95"
$ws.Range("E28").Value = "  +1.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = $refStyle
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "519.05"
$ws.Range("D30").Style = $refStyle
$ws.Range("E30").Value = "  +1.49%  "

$ws.Range("D31").Value = "0.0₃0887"
$ws.Range("E31").Value = "  -2.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.77"
$ws.Range("D32").Style = $refStyle
$ws.Range("E32").Value = "  -0.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.24"
$ws.Range("D33").Style = $refStyle
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.04"
$ws.Range("D36").Style = $refStyle
$ws.Range("E36").Value = "  +1.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.119"
$ws.Range("D37").Style = $refStyle
$ws.Range("E37").Value = "  -2.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.63"
$ws.Range("D38").Style = $refStyle
$ws.Range("E38").Value = "  +1.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.67"
$ws.Range("D39").Style = $refStyle
$ws.Range("E39").Value = "  -0.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("D40").Style = $refStyle
$ws.Range("E40").Value = "  +2.92%  "

$ws.Range("E41").Value = "  -1.32%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.326"
$ws.Range("D43").Style = $refStyle
$ws.Range("E43").Value = "  -1.21%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.82"
$ws.Range("D44").Style = $refStyle
$ws.Range("E44").Value = "  -1.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("D45").Style = $refStyle
$ws.Range("E45").Value = "  -2.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "151.63"
$ws.Range("D46").Style = $refStyle
$ws.Range("E46").Value = "  +5.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.58"
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = "  +2.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.516"
$ws.Range("D48").Style = $refStyle
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0739"
$ws.Range("D49").Style = $refStyle
$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.58"
$ws.Range("D50").Style = $refStyle
$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.577"
$ws.Range("D51").Style = $refStyle
$ws.Range("E51").Value = "  -1.62%  "
